$d = $word.ActiveDocument

# --- Paragraph "5CFD60E0" (empty paragraph before "专业排名相对靠前"):
#     remove the stray <w:rFonts w:hint="eastAsia"/> from the paragraph
#     mark's run properties (pPr/rPr). Nothing else about the paragraph
#     changes, so re-emit it verbatim minus that one element.
$p7 = $d.Paragraphs.Item(7).Range
$xml7 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5CFD60E0" w14:textId="77777777" w:rsidR="001933B3" w:rsidRDefault="001933B3"><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p7.InsertXML($xml7)

# --- Paragraph "7184CC29" ("专业排名相对靠前"):
#     remove the same stray <w:rFonts w:hint="eastAsia"/> from the
#     paragraph mark's run properties (pPr/rPr) while leaving the run
#     of text (which keeps its own rFonts hint) untouched.
$p8 = $d.Paragraphs.Item(8).Range
$xml8 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="7184CC29" w14:textId="4EC8E3C8" w:rsidR="001933B3" w:rsidRPr="000548EB" w:rsidRDefault="001933B3"><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>专业排名相对靠前</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p8.InsertXML($xml8)

# --- Paragraph "78071145" (last, empty paragraph at the end of the
#     document): add a new run with the "个人品质：..." text, matching
#     the formatting used throughout the rest of the resume.
$p9 = $d.Paragraphs.Item(9).Range
$p9inner = $d.Range($p9.Start, $p9.End - 1)
$newText = "个人品质：勤劳、上进、有爱心"
$xml9 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>' + $newText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p9inner.InsertXML($xml9)

Write-Host "Edit applied"
